# Update "想去人数" (want-to-go count) figures for the two sheets that
# contain the full event listing: 展览 (Exhibitions) and 全部类型 (All types).
# Both sheets hold the same four data rows (rows 2-5) in column F.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 656
    $ws.Range("F3").Value = 3898
    $ws.Range("F4").Value = 108
    $ws.Range("F5").Value = 731
}
